# Auto-generated Excel COM-interop script applying Shiva_Profits market-data refresh
# across all 8 crafting-class worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 6667104
$ws.Range("I8").Value = 9091215
$ws.Range("J8").Value = 799.75
$ws.Range("K8").Value = 27273645
$ws.Range("L8").Value = 2399.25
$ws.Range("M8").Value = -27273506
$ws.Range("N8").Value = -2677.25
$ws.Range("H106").Value = 11114660
$ws.Range("I106").Value = 15876077
$ws.Range("K106").Value = 15876077
$ws.Range("M106").Value = -15875446
$ws.Range("H107").Value = 22733104
$ws.Range("I107").Value = 33337800
$ws.Range("J107").Value = 8754.143
$ws.Range("K107").Value = 33337800
$ws.Range("L107").Value = 8754.143
$ws.Range("M107").Value = -33335880
$ws.Range("N107").Value = -12594.143
$ws.Range("H116").Value = 4147.2144
$ws.Range("I116").Value = 5082.7334
$ws.Range("J116").Value = 3627.4814
$ws.Range("K116").Value = 5082.7334
$ws.Range("L116").Value = 3627.4814
$ws.Range("M116").Value = -1640.7334
$ws.Range("N116").Value = -10511.4814
$ws.Range("H132").Value = 8594.139999999999
$ws.Range("I132").Value = 4927.25
$ws.Range("K132").Value = 14781.75
$ws.Range("M132").Value = -12251.75
$ws.Range("H135").Value = 1158.5106
$ws.Range("I135").Value = 789.7907
$ws.Range("K135").Value = 7108.1163
$ws.Range("M135").Value = -4573.1163
$ws.Range("H138").Value = 15026775
$ws.Range("I138").Value = 34484764
$ws.Range("J138").Value = 1591496.4
$ws.Range("K138").Value = 103454292
$ws.Range("L138").Value = 4774489.199999999
$ws.Range("M138").Value = -103449152
$ws.Range("N138").Value = -4784769.199999999
$ws.Range("H140").Value = 68885.39999999999
$ws.Range("J140").Value = 68885.39999999999
$ws.Range("L140").Value = 68885.39999999999
$ws.Range("N140").Value = -79245.39999999999
$ws.Range("H141").Value = 4044.2727
$ws.Range("I141").Value = 4044.2727
$ws.Range("K141").Value = 12132.8181
$ws.Range("M141").Value = -6952.8181

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1995.64
$ws.Range("I32").Value = 1954.2084
$ws.Range("K32").Value = 1954.2084
$ws.Range("M32").Value = -1667.2084
$ws.Range("H97").Value = 52686052
$ws.Range("I97").Value = 58825560
$ws.Range("K97").Value = 58825560
$ws.Range("M97").Value = -58825064
$ws.Range("H102").Value = 1611.3182
$ws.Range("I102").Value = 1529.525
$ws.Range("K102").Value = 1529.525
$ws.Range("M102").Value = 92.47499999999991
$ws.Range("H117").Value = 144000
$ws.Range("J117").Value = 144000
$ws.Range("L117").Value = 144000
$ws.Range("N117").Value = -153178
$ws.Range("H139").Value = 167741.4
$ws.Range("J139").Value = 167741.4
$ws.Range("L139").Value = 167741.4
$ws.Range("N139").Value = -178021.4
$ws.Range("H141").Value = 192496
$ws.Range("J141").Value = 192496
$ws.Range("L141").Value = 192496
$ws.Range("N141").Value = -202856

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3548.973
$ws.Range("I20").Value = 3285.25
$ws.Range("J20").Value = 4369.4443
$ws.Range("K20").Value = 3285.25
$ws.Range("L20").Value = 4369.4443
$ws.Range("M20").Value = -3038.25
$ws.Range("N20").Value = -4863.4443
$ws.Range("H97").Value = 2766.4
$ws.Range("I97").Value = 2766.4
$ws.Range("K97").Value = 2766.4
$ws.Range("M97").Value = -1775.4
$ws.Range("H99").Value = 1581.2354
$ws.Range("I99").Value = 1418.5
$ws.Range("J99").Value = 2340.6667
$ws.Range("K99").Value = 1418.5
$ws.Range("L99").Value = 2340.6667
$ws.Range("M99").Value = 79.5
$ws.Range("N99").Value = -5336.6667
$ws.Range("H105").Value = 3299.8635
$ws.Range("I105").Value = 1425.0526
$ws.Range("J105").Value = 15173.667
$ws.Range("K105").Value = 1425.0526
$ws.Range("L105").Value = 15173.667
$ws.Range("M105").Value = 321.9474
$ws.Range("N105").Value = -18667.667
$ws.Range("H120").Value = 112880.5
$ws.Range("J120").Value = 112880.5
$ws.Range("L120").Value = 112880.5
$ws.Range("N120").Value = -122556.5
$ws.Range("H125").Value = 49999
$ws.Range("J125").Value = 49999
$ws.Range("L125").Value = 49999
$ws.Range("N125").Value = -59839

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 557.5625
$ws.Range("I22").Value = 422.2
$ws.Range("J22").Value = 783.1667
$ws.Range("K22").Value = 422.2
$ws.Range("L22").Value = 783.1667
$ws.Range("M22").Value = -72.19999999999999
$ws.Range("N22").Value = -1483.1667
$ws.Range("H58").Value = 1253.8085
$ws.Range("I58").Value = 1190.3658
$ws.Range("K58").Value = 1190.3658
$ws.Range("M58").Value = -987.3658
$ws.Range("H110").Value = 149980.5
$ws.Range("J110").Value = 149980.5
$ws.Range("L110").Value = 149980.5
$ws.Range("N110").Value = -158160.5
$ws.Range("H124").Value = 59734.332
$ws.Range("J124").Value = 59734.332
$ws.Range("L124").Value = 59734.332
$ws.Range("N124").Value = -64644.332
$ws.Range("H132").Value = 313064.56
$ws.Range("I132").Value = 11075.348
$ws.Range("K132").Value = 33226.044
$ws.Range("M132").Value = -30696.044
$ws.Range("H136").Value = 1253.8085
$ws.Range("I136").Value = 1190.3658
$ws.Range("K136").Value = 3571.0974
$ws.Range("M136").Value = -1021.0974

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 5064
$ws.Range("I116").Value = 5876.3335
$ws.Range("J116").Value = 2627
$ws.Range("K116").Value = 17629.0005
$ws.Range("L116").Value = 7881
$ws.Range("M116").Value = -14187.0005
$ws.Range("N116").Value = -14765
$ws.Range("H131").Value = 1471344.9
$ws.Range("J131").Value = 1599.75
$ws.Range("L131").Value = 4799.25
$ws.Range("N131").Value = -14879.25
$ws.Range("H134").Value = 1427.7255
$ws.Range("I134").Value = 1123.7021
$ws.Range("K134").Value = 3371.1063
$ws.Range("M134").Value = 1698.8937
$ws.Range("H137").Value = 4443.8823
$ws.Range("J137").Value = 3942.9167
$ws.Range("L137").Value = 11828.7501
$ws.Range("N137").Value = -22028.7501

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 30003862
$ws.Range("I80").Value = 68573980
$ws.Range("K80").Value = 68573980
$ws.Range("M80").Value = -68572982
$ws.Range("H83").Value = 30003862
$ws.Range("I83").Value = 68573980
$ws.Range("K83").Value = 342869900
$ws.Range("M83").Value = -342864908
$ws.Range("H97").Value = 1566.7646
$ws.Range("I97").Value = 1505.4
$ws.Range("J97").Value = 1654.4286
$ws.Range("K97").Value = 1505.4
$ws.Range("L97").Value = 1654.4286
$ws.Range("M97").Value = -1009.4
$ws.Range("N97").Value = -2646.4286
$ws.Range("H141").Value = 44319.332
$ws.Range("J141").Value = 44319.332
$ws.Range("L141").Value = 44319.332
$ws.Range("N141").Value = -54679.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 989
$ws.Range("I22").Value = 988.125
$ws.Range("J22").Value = 991.3333
$ws.Range("K22").Value = 988.125
$ws.Range("L22").Value = 991.3333
$ws.Range("M22").Value = -693.125
$ws.Range("N22").Value = -1581.3333
$ws.Range("H27").Value = 989
$ws.Range("I27").Value = 988.125
$ws.Range("J27").Value = 991.3333
$ws.Range("K27").Value = 988.125
$ws.Range("L27").Value = 991.3333
$ws.Range("M27").Value = -881.125
$ws.Range("N27").Value = -1205.3333
$ws.Range("H46").Value = 2646.3635
$ws.Range("I46").Value = 922.9
$ws.Range("K46").Value = 922.9
$ws.Range("M46").Value = -734.9
$ws.Range("H130").Value = 150000
$ws.Range("J130").Value = 150000
$ws.Range("L130").Value = 150000
$ws.Range("N130").Value = -160040
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 96463.336
$ws.Range("J46").Value = 120000
$ws.Range("L46").Value = 120000
$ws.Range("N46").Value = -120462
$ws.Range("H96").Value = 4514.75
$ws.Range("I96").Value = 4793.25
$ws.Range("K96").Value = 4793.25
$ws.Range("M96").Value = -3420.25
$ws.Range("H100").Value = 997.1429000000001
$ws.Range("I100").Value = 1230
$ws.Range("J100").Value = 822.5
$ws.Range("K100").Value = 2460
$ws.Range("L100").Value = 1645
$ws.Range("M100").Value = -1919
$ws.Range("N100").Value = -2727
$ws.Range("H105").Value = 64800
$ws.Range("J105").Value = 64800
$ws.Range("L105").Value = 64800
$ws.Range("N105").Value = -71788
$ws.Range("H107").Value = 1782.6666
$ws.Range("I107").Value = 1139.2
$ws.Range("K107").Value = 3417.6
$ws.Range("M107").Value = -1497.6
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H121").Value = 167999
$ws.Range("J121").Value = 167999
$ws.Range("L121").Value = 167999
$ws.Range("N121").Value = -171493
$ws.Range("H132").Value = 4919.2905
$ws.Range("I132").Value = 4274.263
$ws.Range("K132").Value = 12822.789
$ws.Range("M132").Value = -10292.789
$ws.Range("H134").Value = 96463.336
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 360000
$ws.Range("N134").Value = -365070
